$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 4668
$ws.Range("B2").Value = "Hellena Carvalho"
$ws.Range("C2").Value = "Engenharia"
$ws.Range("D2").Value = "Consulta medica"
$ws.Range("E2").Value = 5
$ws.Range("F2").Value = 45086
$ws.Range("G2").Value = 4366.22

# Row 3
$ws.Range("A3").Value = 30912
$ws.Range("B3").Value = "Maria Cecília Carvalho"
$ws.Range("C3").Value = "Financeiro"
$ws.Range("D3").Value = "Doenca"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 45099
$ws.Range("G3").Value = 5952.12

# Row 4
$ws.Range("A4").Value = 69195
$ws.Range("B4").Value = "Luísa da Cruz"
$ws.Range("C4").Value = "Engenharia"
$ws.Range("D4").Value = "Consulta medica"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 45081
$ws.Range("G4").Value = 2590.57

# Row 5
$ws.Range("A5").Value = 57730
$ws.Range("B5").Value = "Dra. Ísis Pereira"
$ws.Range("C5").Value = "Juridico"
$ws.Range("D5").Value = "Viagem de negocios"
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = 45102
$ws.Range("G5").Value = 3936.32

# Row 6
$ws.Range("A6").Value = 65903
$ws.Range("B6").Value = "Ravi Pastor"
$ws.Range("C6").Value = "P&D"
$ws.Range("D6").Value = "Outros"
$ws.Range("E6").Value = 8
$ws.Range("F6").Value = 45086
$ws.Range("G6").Value = 6122.55

# Row 7
$ws.Range("A7").Value = 13756
$ws.Range("B7").Value = "Rafael Lima"
$ws.Range("C7").Value = "Vendas"
$ws.Range("D7").Value = "Viagem de negocios"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 45100
$ws.Range("G7").Value = 2334.19

# Row 8
$ws.Range("A8").Value = 24931
$ws.Range("B8").Value = "Liam Farias"
$ws.Range("C8").Value = "P&D"
$ws.Range("D8").Value = "Doenca"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 45087
$ws.Range("G8").Value = 8782.02

# Row 9
$ws.Range("A9").Value = 27939
$ws.Range("B9").Value = "Sra. Alícia Viana"
$ws.Range("C9").Value = "Marketing"
$ws.Range("D9").Value = "Doenca"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 45088
$ws.Range("G9").Value = 2749.05

# Row 10
$ws.Range("A10").Value = 5464
$ws.Range("B10").Value = "Sra. Ana Cecília Pimenta"
$ws.Range("C10").Value = "Financeiro"
$ws.Range("D10").Value = "Consulta medica"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 45094
$ws.Range("G10").Value = 9549.940000000001

# Row 11
$ws.Range("A11").Value = 62744
$ws.Range("B11").Value = "Ana Lívia Aparecida"
$ws.Range("C11").Value = "TI"
$ws.Range("D11").Value = "Viagem de negocios"
$ws.Range("E11").Value = 7
$ws.Range("F11").Value = 45086
$ws.Range("G11").Value = 5539.58
